$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.182.30'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '1.563.41'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  +0.09%  '
$c = $ws.Range('D5')
$c.Value = "'206.86"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('E6').Value = '  -1.86%  '
$ws.Range('E7').Value = '  +0.10%  '
$c = $ws.Range('D8')
$c.Value = "'22.05"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.96%  '
$c = $ws.Range('D9')
$c.Value = "'0.247"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').Value = '1.784.54'
$ws.Range('E12').Value = '  -1.72%  '
$ws.Range('D13').Value = '1.563.58'
$ws.Range('E13').Value = '  -1.38%  '
$c = $ws.Range('D14')
$c.Value = "'3.76"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.68%  '
$c = $ws.Range('D15')
$c.Value = "'0.516"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -2.76%  '
$c = $ws.Range('D16')
$c.Value = "'63.07"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.78%  '
$ws.Range('D17').Value = '27.200.39'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').Value = '0.0₃0688'
$ws.Range('E18').Value = '  -1.19%  '
$c = $ws.Range('D19')
$c.Value = "'211.91"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -3.90%  '
$c = $ws.Range('D20')
$c.Value = "'7.22"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -1.86%  '
$ws.Range('E21').Value = '  +0.10%  '
$c = $ws.Range('D22')
$c.Value = "'4.10"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.05%  '
$c = $ws.Range('D23')
$c.Value = "'9.43"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.81%  '
$c = $ws.Range('D24')
$c.Value = "'1.98"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.47%  '
$c = $ws.Range('D25')
$c.Value = "'152.40"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  -3.69%  '
$c = $ws.Range('D27')
$c.Value = "'14.83"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -1.80%  '
$ws.Range('E30').Value = '  -0.72%  '
$c = $ws.Range('D31')
$c.Value = "'0.0464"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('D33').Value = '1.374.85'
$ws.Range('E33').Value = '  +0.15%  '
$c = $ws.Range('D34')
$c.Value = "'2.94"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.30%  '
$c = $ws.Range('D35')
$c.Value = "'1.55"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.51%  '
$c = $ws.Range('D37')
$c.Value = "'0.942"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -3.86%  '
$ws.Range('E38').Value = '  -1.57%  '
$c = $ws.Range('D39')
$c.Value = "'0.522"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -3.38%  '
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range('E41').Value = '  +0.12%  '
$c = $ws.Range('D42')
$c.Value = "'0.989"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +2.06%  '
$c = $ws.Range('D43')
$c.Value = "'1.79"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +2.97%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D44')
$c.Value = "'63.42"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -1.81%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D45')
$c.Value = "'2.17"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.29%  '
$c = $ws.Range('D46')
$c.Value = "'5.21"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').Value = '1.697.73'
$ws.Range('E47').Value = '  -1.55%  '
$c = $ws.Range('D48')
$c.Value = "'85.49"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.73%  '
$ws.Range('D49').Value = '0.0₇0995'
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$c = $ws.Range('D51')
$c.Value = "'1.00"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.23%  '
